{"js": "// Update the date paragraph and the 25 \"three-digit \u00d7 one-digit\"\n// multiplication prompts laid out in a 5-column table (5 rows of\n// problems, each followed by 4 blank spacer rows).\n//\n// Every replacement is done positionally (by paragraph index / table\n// row+column index) rather than by text search, because one of the new\n// values (\"877\u00d75=\") is identical to one of the *old* values elsewhere in\n// the document \u2014 a global find/replace could clobber the wrong cell.\n\nconst body = context.document.body;\n\n// --- 1) Date line: first paragraph of the body -----------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2024-10-09 Wednesday\", \"Replace\");\n\n// --- 2) Multiplication problems in the table --------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, in row-major order, for the 5 table rows that hold\n// problems (rows 0, 4, 9, 14, 19 of the 20-row table).\nconst newGrid = [\n  [\"877\u00d75=\", \"622\u00d78=\", \"985\u00d76=\", \"269\u00d78=\", \"971\u00d79=\"],\n  [\"220\u00d72=\", \"143\u00d76=\", \"456\u00d72=\", \"293\u00d75=\", \"175\u00d74=\"],\n  [\"563\u00d76=\", \"481\u00d78=\", \"849\u00d78=\", \"728\u00d73=\", \"244\u00d79=\"],\n  [\"515\u00d78=\", \"139\u00d72=\", \"108\u00d76=\", \"118\u00d73=\", \"899\u00d72=\"],\n  [\"331\u00d75=\", \"541\u00d72=\", \"561\u00d79=\", \"447\u00d72=\", \"586\u00d78=\"],\n];\nconst tableRows = [0, 4, 9, 14, 19];\n\nfor (let r = 0; r < tableRows.length; r++) {\n  const rowIndex = tableRows[r];\n  for (let c = 0; c < newGrid[r].length; c++) {\n    const cell = table.getCell(rowIndex, c);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    cellParagraphs.items[0].insertText(newGrid[r][c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the 25 \"three-digit x one-digit\"\n# multiplication prompts laid out in a 5-column table (5 rows of\n# problems, each followed by 4 blank spacer rows).\n#\n# Every replacement is done positionally (by paragraph index / table\n# row+column index) rather than by text search, because one of the new\n# values (\"877x5=\") is identical to one of the *old* values elsewhere in\n# the document - a global find/replace could clobber the wrong cell.\n# Setting Range.Text (rather than rebuilding the run) keeps each run's\n# existing formatting (font, size, alignment) untouched.\n\n$d = $word.ActiveDocument\n\n# --- 1) Date line: first paragraph of the body -----------------------\n$d.Paragraphs.Item(1).Range.Text = \"2024-10-09 Wednesday\"\n\n# --- 2) Multiplication problems in the table --------------------------\n$table = $d.Tables.Item(1)\n\n# New values, in row-major order, for the 5 table rows that hold\n# problems (1-based rows 1, 5, 10, 15, 20 of the 20-row table).\n$newGrid = @(\n    @(\"877\u00d75=\", \"622\u00d78=\", \"985\u00d76=\", \"269\u00d78=\", \"971\u00d79=\"),\n    @(\"220\u00d72=\", \"143\u00d76=\", \"456\u00d72=\", \"293\u00d75=\", \"175\u00d74=\"),\n    @(\"563\u00d76=\", \"481\u00d78=\", \"849\u00d78=\", \"728\u00d73=\", \"244\u00d79=\"),\n    @(\"515\u00d78=\", \"139\u00d72=\", \"108\u00d76=\", \"118\u00d73=\", \"899\u00d72=\"),\n    @(\"331\u00d75=\", \"541\u00d72=\", \"561\u00d79=\", \"447\u00d72=\", \"586\u00d78=\")\n)\n$tableRows = @(1, 5, 10, 15, 20)\n\nfor ($r = 0; $r -lt $tableRows.Length; $r++) {\n    $rowIndex = $tableRows[$r]\n    $rowValues = $newGrid[$r]\n    for ($c = 0; $c -lt $rowValues.Length; $c++) {\n        $cell = $table.Cell($rowIndex, $c + 1)\n        $cell.Range.Text = $rowValues[$c]\n    }\n}\n"}
